$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2021-12-08"

# Update row 14 (December) label
$ws.Range("A14").Value = "December (through 12-08)"

# 2015 columns (B,C,D = arrest_made, no_arrest_made, arrest_rate)
$ws.Range("C14").Value = 6
$ws.Range("D14").Value = 0.25

# 2016 columns (E,F,G)
$ws.Range("E14").Value = 2
$ws.Range("G14").Value = 0.08699999999999999

# 2017 columns (H,I,J)
$ws.Range("I14").Value = 27
$ws.Range("J14").Value = 0.1

# 2018 columns (K,L,M)
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 15
$ws.Range("M14").Value = 0.1176

# 2019 columns (N,O,P)
$ws.Range("N14").Value = 3
$ws.Range("P14").Value = 0.3

# 2020 columns (Q,R,S)
$ws.Range("R14").Value = 37
$ws.Range("S14").Value = 0.0513

# 2021 columns (T,U,V) - only U14 present
$ws.Range("U14").Value = 65

# Row 15 (Total)
$ws.Range("C15").Value = 264
$ws.Range("D15").Value = 0.1171

$ws.Range("E15").Value = 62
$ws.Range("G15").Value = 0.1058

$ws.Range("I15").Value = 785
$ws.Range("J15").Value = 0.0776

$ws.Range("K15").Value = 76
$ws.Range("L15").Value = 623
$ws.Range("M15").Value = 0.1087

$ws.Range("N15").Value = 57
$ws.Range("P15").Value = 0.1048

$ws.Range("R15").Value = 1237
$ws.Range("S15").Value = 0.0507

$ws.Range("U15").Value = 1608
$ws.Range("V15").Value = 0.058
